$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.08921075830798841
$ws.Range("D3").Value = 0.08469939697638122
$ws.Range("D4").Value = 0.09011205197555508
$ws.Range("D5").Value = 0.07624338743113865
$ws.Range("D6").Value = 0.06989478692934439
$ws.Range("D7").Value = 0.06371476635184599
$ws.Range("D8").Value = 0.1256107809981147
$ws.Range("D9").Value = 0.1211262762080012
$ws.Range("D10").Value = 0.110924843000073
$ws.Range("D11").Value = 0.09876737288663349
$ws.Range("D12").Value = 0.0876510592858345
$ws.Range("D13").Value = 0.07850770639489048
